$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Logout"
$ws.Range("B5").Value = "noch nicht auf allen Seiten implementiert"

$ws.Range("B5").Select()
